$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet4 -> Sheet7)
$ws.Name = "Sheet7"

# Update the process name: "Induction Hardening Bearing Surface 1" ->
# "Induction Hardening Bearing Surfaces 1, 2" (less impactful hardening process
# now covers both bearing surfaces)
$ws.Range("X17").Value = "Induction Hardening Bearing Surfaces 1, 2"
$ws.Range("B21").Value = "Induction Hardening Bearing Surfaces 1, 2"

# Update run date/time (D1, F1) - use raw Excel serial values so the
# existing date/time number formats are preserved
$ws.Range("D1").Value = 45572
$ws.Range("F1").Value = 0.810595706018518

# Update the impact values for Induction Hardening (reduced impact) and the
# corresponding overall totals
$ws.Range("E18").Value = 173.688803990365
$ws.Range("X18").Value = 34.074702596165
$ws.Range("E21").Value = 34.074702596165
$ws.Range("X21").Value = 34.074702596165
